# Weekly crypto price/volume refresh (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to keep a literal text value even when the
    # string looks numeric (e.g. "0.999"), matching the source
    # data which stores every Price/Volume entry as text, then
    # restore the default (unstyled) look of the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "59.103.31"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "2.494.96"
$ws.Range("E3").Value = "  -1.38%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "0.999"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "534.34"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "136.67"
$ws.Range("E6").Value = "  -2.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.563"
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "2.517.31"
$ws.Range("E9").Value = "  -0.80%  "

# Row 10
$ws.Range("E10").Value = "  +1.56%  "

# Row 11
$ws.Range("E11").Value = "  -0.45%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "5.35"
$ws.Range("E12").Value = "  -0.67%  "

# Row 13
$ws.Range("E13").Value = "  -2.52%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "2.942.29"
$ws.Range("E14").Value = "  -1.21%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "23.20"
$ws.Range("E15").Value = "  +0.23%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "58.822.48"
$ws.Range("E16").Value = "  -0.92%  "

# Row 17
$ws.Range("E17").Value = "  -1.00%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "2.515.54"
$ws.Range("E18").Value = "  -0.06%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "11.05"
$ws.Range("E19").Value = "  +0.72%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "4.25"
$ws.Range("E20").Value = "  +0.50%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "324.41"
$ws.Range("E21").Value = "  +0.67%  "

# Row 22
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "5.86"
$ws.Range("E23").Value = "  +0.57%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "64.48"
$ws.Range("E24").Value = "  +4.16%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "0.419"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26
$ws.Range("E26").Value = "  -0.53%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "0.998"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "7.57"
$ws.Range("E28").Value = "  -3.25%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "6.75"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "0.0₃0774"
$ws.Range("E30").Value = "  +0.42%  "

# Row 31
$ws.Range("E31").Value = "  -1.99%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "167.87"
$ws.Range("E32").Value = "  +4.02%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Cells.Item(33, 4) "1.17"
$ws.Range("E33").Value = "  +3.36%  "

# Row 34
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Cells.Item(34, 4) "0.998"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "1.42"
$ws.Range("E35").Value = "  -3.00%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "18.53"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "4.09"
$ws.Range("E37").Value = "  -3.16%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "1.57"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "36.75"
$ws.Range("E39").Value = "  -0.82%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.825"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "3.62"
$ws.Range("E41").Value = "  -0.98%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "5.26"
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "279.07"
$ws.Range("E43").Value = "  -2.49%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.996"
$ws.Range("E44").Value = "  -0.17%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "0.603"
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "10.89"
$ws.Range("E46").Value = "  +0.27%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "128.81"
$ws.Range("E47").Value = "  +3.68%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "0.0930"
$ws.Range("E48").Value = "  +0.49%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "0.0515"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "0.0221"
$ws.Range("E50").Value = "  -0.86%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "17.37"
$ws.Range("E51").Value = "  -1.22%  "
